$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'312.52"
$ws.Range("E2").Value = "'8.74%"
$ws.Range("D3").Value = "'32.42"
$ws.Range("E3").Value = "'9.48%"
$ws.Range("D4").Value = "'5.338"
$ws.Range("E4").Value = "'4.44%"
$ws.Range("D5").Value = "'0.07674"
$ws.Range("E5").Value = "'14.63%"
$ws.Range("D6").Value = "'7.875"
$ws.Range("E6").Value = "'7.34%"
$ws.Range("D7").Value = "'3.714"
$ws.Range("E7").Value = "'9.17%"
$ws.Range("D8").Value = "'1.602"
$ws.Range("E8").Value = "'18.00%"
$ws.Range("D9").Value = "'0.9198"
$ws.Range("E9").Value = "'0.81%"
$ws.Range("D10").Value = "'0.01714"
$ws.Range("E10").Value = "'2,555.28%"
$ws.Range("D11").Value = "'0.1728"
$ws.Range("E11").Value = "'8.85%"
$ws.Range("D12").Value = "'0.07574"
$ws.Range("E12").Value = "'12.27%"
$ws.Range("D13").Value = "'0.08244"
$ws.Range("E13").Value = "'7.03%"
$ws.Range("D14").Value = "'0.03028"
$ws.Range("E14").Value = "'3.27%"
$ws.Range("D15").Value = "'0.09901"
$ws.Range("E15").Value = "'10.30%"
$ws.Range("D16").Value = "'0.001521"
$ws.Range("E16").Value = "'-3.43%"
$ws.Range("D17").Value = "'0.04566"
$ws.Range("E17").Value = "'1.72%"
$ws.Range("D18").Value = "'0.006268"
$ws.Range("E18").Value = "'-0.09%"
$ws.Range("D19").Value = "'3.480"
$ws.Range("E19").Value = "'1.20%"
$ws.Range("E20").Value = "'0.94%"
$ws.Range("E21").Value = "'3.16%"
$ws.Range("E22").Value = "'0.72%"
$ws.Range("D23").Value = "'4.250"
$ws.Range("E23").Value = "'4.54%"
$ws.Range("D24").Value = "'0.1626"
$ws.Range("E24").Value = "'2.80%"
$ws.Range("D25").Value = "'0.001220"
$ws.Range("E25").Value = "'2.42%"
$ws.Range("D26").Value = "'0.004501"
$ws.Range("E26").Value = "'9.37%"
$ws.Range("E27").Value = "'8.15%"
$ws.Range("D28").Value = "'0.0001740"
$ws.Range("E28").Value = "'7.65%"
$ws.Range("D40").Value = "'0.04642"
$ws.Range("E40").Value = "'8.81%"
$ws.Range("D41").Value = "'0.007217"
$ws.Range("E41").Value = "'7.29%"
$ws.Range("D42").Value = "'0.1374"
$ws.Range("E42").Value = "'10.81%"
$ws.Range("D43").Value = "'0.002255"
$ws.Range("E43").Value = "'1.17%"
$ws.Range("D44").Value = "'0.01432"
$ws.Range("E44").Value = "'7.20%"
$ws.Range("D45").Value = "'0.00006186"
$ws.Range("E45").Value = "'9.16%"
$ws.Range("D46").Value = "'1.892"
$ws.Range("E46").Value = "'-4.13%"
$ws.Range("E47").Value = "'-0.49%"
